# Implemented PUT with status transitions. Also CORs for the Web API.
#
# The "status transition" mini-table (columns G:K, rows 5-9) shifts by one
# step: a new "Cancelled" transition is inserted at row 8 (col H) and the
# old one at row 9 (col H) is removed, with every following transition in
# that chain sliding up/left by one column. The unused "InProgress" shared
# string is dropped automatically once no cell references it any more.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yellow = 65535  # RGB(255,255,0) as an OLE COLORREF (R + G*256 + B*65536)

# --- Remove the cells that disappear entirely in the new layout ---------
$ws.Range("J7").Clear()
$ws.Range("K7").Clear()
$ws.Range("H9").Clear()

# --- Row 5: WIP -> {Complete, Archive} instead of {InProgress, WIP} -----
$ws.Range("I5").Value = "Complete"
$ws.Range("J5").Value = "Archive"
$ws.Range("J5").Interior.Color = $yellow

# --- Row 6: Cancelled -> {FollowUp, Complete, Archive} -------------------
$ws.Range("I6").Value = "FollowUp"
$ws.Range("J6").Value = "Complete"
$ws.Range("J6").Style = "Normal"
$ws.Range("K6").Value = "Archive"
$ws.Range("K6").Interior.Color = $yellow

# --- Row 7: InProgress -> {Backlog} --------------------------------------
$ws.Range("I7").Value = "Backlog"

# --- Row 8: Complete -> {Cancelled, Backlog, <blank>} --------------------
$ws.Range("H8").Value = "Cancelled"
$ws.Range("H8").Interior.Color = $yellow

# --- Row 9: FollowUp no longer has an outgoing "Cancelled" transition ----
# (handled by the Clear() above)

# --- Move the active selection/scroll the way the author left the sheet -
$win = $excel.ActiveWindow
$ws.Range("H9").Select()
$win.ScrollColumn = 5
$win.ScrollRow = 1
